$d = $word.ActiveDocument

function Get-ContainingParagraphRange($pos) {
    foreach ($pp in $d.Content.Paragraphs) {
        $s = $pp.Range.Start
        $e = $pp.Range.End
        if ($pos -ge $s -and $pos -lt $e) {
            return $pp.Range
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Hunk 1: "Meeting Topic" table cell - "Kick Off Meeting" -> "Status- und
# Planungsmeeting" split across three runs, with the (moved) "_GoBack"
# bookmark sitting between the "-" run and the " und Planungsmeeting" run.
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Kick Off Meeting")
if (-not $found1) {
    throw "Could not find 'Kick Off Meeting'"
}
$para1 = Get-ContainingParagraphRange($rng1.Start)
if ($para1 -eq $null) {
    throw "Could not locate containing paragraph for 'Kick Off Meeting'"
}

$xml1 = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="007662FB" w:rsidRPr="00D57F14" w:rsidRDefault="009B7DA6" w:rsidP="00782A7D"><w:pPr><w:rPr><w:color w:val="333333"/><w:lang w:val="de-AT"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="333333"/><w:lang w:val="de-AT"/></w:rPr><w:t>Status</w:t></w:r><w:r><w:rPr><w:color w:val="333333"/><w:lang w:val="de-AT"/></w:rPr><w:t>-</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:color w:val="333333"/><w:lang w:val="de-AT"/></w:rPr><w:t xml:space="preserve"> und Planungsmeeting</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$para1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Hunk 2: "Action Items" / sync paragraph - drop the old "_GoBack" bookmark
# that used to sit right after "abgleichen" (the bookmark moved to hunk 1).
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("abgleichen")
if (-not $found2) {
    throw "Could not find 'abgleichen'"
}
$para2 = Get-ContainingParagraphRange($rng2.Start)
if ($para2 -eq $null) {
    throw "Could not locate containing paragraph for 'abgleichen'"
}

$xml2 = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="007662FB" w:rsidRPr="00682CBF" w:rsidRDefault="00BF4F6D" w:rsidP="00BF4F6D"><w:pPr><w:tabs><w:tab w:val="left" w:pos="1605"/></w:tabs><w:rPr><w:color w:val="333333"/><w:lang w:val="de-AT"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Projektplan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>zeitnahe</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>überlegen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> und </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>abgleichen</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$para2.InsertXML($xml2)

Write-Host "Done."
